$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "Davangere"
$ws.Range("G5").Value = "Davangere"
$ws.Range("G6").Value = "Chikballapur"
$ws.Range("G7").Value = "Davangere"
$ws.Range("G8").Value = "Davangere"
$ws.Range("G9").Value = "Davangere"
$ws.Range("G11").Value = "Davangere"
$ws.Range("G12").Value = "Chikballapur"
$ws.Range("G16").Value = "Bagalkot"
$ws.Range("G17").Value = "Bagalkot"
$ws.Range("G18").Value = "Bagalkot"
$ws.Range("F19").ClearContents()
$ws.Range("G22").Value = "Chikballapur"
$ws.Range("G23").Value = "Chikballapur"
$ws.Range("G24").Value = "Bagalkot"
$ws.Range("G25").Value = "Chikballapur"
$ws.Range("G27").Value = "Davangere"
$ws.Range("G28").Value = "Chikballapur"
$ws.Range("G30").Value = "Davangere"
$ws.Range("G31").Value = "Bagalkot"
$ws.Range("G33").Value = "Chikballapur"
$ws.Range("G34").Value = "Chikballapur"
$ws.Range("G35").Value = "Vijayapura (Bijapur)"
$ws.Range("G36").Value = "Davangere"
$ws.Range("G37").Value = "Davangere"
$ws.Range("G38").Value = "Davangere"
$ws.Range("G40").Value = "Davangere"
$ws.Range("G41").Value = "Chikballapur"
$ws.Range("G42").Value = "Davangere"
$ws.Range("G43").Value = "Bagalkot"
$ws.Range("G45").Value = "Davangere"
$ws.Range("G46").Value = "Davangere"
$ws.Range("G47").Value = "Chikballapur"
$ws.Range("G48").Value = "Chikballapur"
$ws.Range("G49").Value = "Chikballapur"
$ws.Range("G50").Value = "Chikballapur"
$ws.Range("F51").ClearContents()
$ws.Range("G51").Value = "Chikballapur"
$ws.Range("G52").Value = "Chikballapur"
$ws.Range("G54").Value = "Bagalkot"
$ws.Range("G58").Value = "Chikballapur"
$ws.Range("G59").Value = "Chikballapur"
$ws.Range("G61").Value = "Bagalkot"
$ws.Range("G63").Value = "Chikballapur"
$ws.Range("G64").Value = "Bagalkot"
$ws.Range("G65").Value = "Chikballapur"
$ws.Range("G66").Value = "Chikballapur"
$ws.Range("F68").ClearContents()
$ws.Range("G68").Value = "Chikballapur"
$ws.Range("F69").ClearContents()
$ws.Range("G69").Value = "Chikballapur"
$ws.Range("G70").Value = "Chikballapur"
$ws.Range("G71").Value = "Chikballapur"
$ws.Range("G73").Value = "Chikballapur"
$ws.Range("G74").Value = "Davangere"
